$wb = $excel.ActiveWorkbook
Write-Output "done"
